$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = "new price (or $null if unchanged)"; E = "new volume" }
$updates = @{
  2 = @{ D = "26.977.41"; E = "  -0.06%  " }
  3 = @{ D = "1.844.56"; E = "  -0.12%  " }
  4 = @{ D = "1.014"; E = "  +0.54%  " }
  5 = @{ D = "1.012"; E = "  +0.46%  " }
  6 = @{ D = "309.21"; E = "  -0.49%  " }
  7 = @{ D = "0.4764"; E = "  +1.89%  " }
  8 = @{ D = $null; E = "  +1.11%  " }
  9 = @{ D = "0.07220"; E = "  +0.52%  " }
  10 = @{ D = "0.9306"; E = "  -0.40%  " }
  11 = @{ D = $null; E = "  +0.95%  " }
  12 = @{ D = "0.07734"; E = "  +0.73%  " }
  13 = @{ D = "1.893.88"; E = "  +3.63%  " }
  14 = @{ D = "5.395"; E = "  +1.78%  " }
  15 = @{ D = "6.454"; E = "  +0.72%  " }
  16 = @{ D = "88.82"; E = "  +0.69%  " }
  17 = @{ D = "1.015"; E = "  +0.61%  " }
  18 = @{ D = "0.000008662"; E = "  +0.84%  " }
  19 = @{ D = $null; E = "  +0.45%  " }
  20 = @{ D = "27.019.37"; E = "  +0.04%  " }
  21 = @{ D = "14.55"; E = "  +1.09%  " }
  22 = @{ D = "5.068"; E = "  +0.61%  " }
  23 = @{ D = "10.64"; E = "  -0.14%  " }
  24 = @{ D = "1.955"; E = "  +1.16%  " }
  25 = @{ D = "152.93"; E = "  +0.09%  " }
  26 = @{ D = "18.22"; E = "  +0.99%  " }
  27 = @{ D = "2.012"; E = "  -0.74%  " }
  28 = @{ D = "114.48"; E = "  +0.43%  " }
  29 = @{ D = "4.967"; E = "  +0.57%  " }
  30 = @{ D = "0.08872"; E = "  +0.24%  " }
  31 = @{ D = "3.312"; E = "  +3.88%  " }
  32 = @{ D = "1.179"; E = "  -0.20%  " }
  33 = @{ D = "0.7427"; E = "  -0.62%  " }
  34 = @{ D = "4.499"; E = "  +0.49%  " }
  35 = @{ D = $null; E = "  -5.14%  " }
  36 = @{ D = "1.118"; E = "  +2.47%  " }
  37 = @{ D = "0.01960"; E = "  +1.04%  " }
  38 = @{ D = "0.05257"; E = "  +1.68%  " }
  39 = @{ D = "2.969"; E = "  -0.62%  " }
  40 = @{ D = "0.5262"; E = "  +2.44%  " }
  41 = @{ D = "7.019"; E = "  +1.55%  " }
  42 = @{ D = $null; E = "  -0.24%  " }
  43 = @{ D = $null; E = "  +1.04%  " }
  44 = @{ D = "10.60"; E = "  +2.14%  " }
  45 = @{ D = "0.4737"; E = "  +0.31%  " }
  46 = @{ D = "1.014"; E = "  +0.48%  " }
  47 = @{ D = "101.95"; E = "  +1.59%  " }
  48 = @{ D = "1.605"; E = "  -0.02%  " }
  49 = @{ D = "65.93"; E = "  +2.70%  " }
  50 = @{ D = "0.06082"; E = "  +0.53%  " }
  51 = @{ D = "0.8912"; E = "  +3.51%  " }
}

foreach ($row in $updates.Keys) {
  $vals = $updates[$row]
  if ($null -ne $vals.D) {
    $cell = $ws.Range("D$row")
    # Force text so numeric-looking strings (single decimal point, e.g.
    # "0.07220") keep their exact formatting (trailing zeros, fixed
    # decimals) instead of being auto-coerced to a number by Excel.
    # Values like "26.977.41" (two dots) are never auto-coerced, so they
    # don't need this treatment.
    if ($vals.D -match '^[+-]?[0-9]*\.?[0-9]+$') {
      $cell.NumberFormat = "@"
      $cell.Value = $vals.D
      $cell.Style = "Normal"
    } else {
      $cell.Value = $vals.D
    }
  }
  if ($null -ne $vals.E) {
    $ws.Range("E$row").Value = $vals.E
  }
}
